# Verify_40V_Calculation_For_PFI_with_Devices.xlsx
# "Updated test data as per new implemenation"
#
# The "Loading Details Name" column on the "Add Panels" sheet lists a
# 40V loading row whose label changes from "40V (A)" to "40V Rail(A)".
# That label is repeated in cells I8, I9 and I10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

$ws.Range("I8").Value  = "40V Rail(A)"
$ws.Range("I9").Value  = "40V Rail(A)"
$ws.Range("I10").Value = "40V Rail(A)"
